{"js": "// Replace each division-problem placeholder text in the document's\n// table cells with its updated value, per the commit diff.\nconst replacements = [\n    { oldText: \"963\u00f77=\", newText: \"708\u00f76=\" },\n    { oldText: \"496\u00f74=\", newText: \"555\u00f76=\" },\n    { oldText: \"877\u00f79=\", newText: \"732\u00f75=\" },\n    { oldText: \"708\u00f73=\", newText: \"257\u00f73=\" },\n    { oldText: \"316\u00f74=\", newText: \"463\u00f72=\" },\n    { oldText: \"314\u00f78=\", newText: \"558\u00f77=\" },\n    { oldText: \"855\u00f74=\", newText: \"233\u00f73=\" },\n    { oldText: \"703\u00f73=\", newText: \"705\u00f77=\" },\n    { oldText: \"128\u00f75=\", newText: \"748\u00f72=\" },\n    { oldText: \"634\u00f74=\", newText: \"694\u00f75=\" },\n    { oldText: \"971\u00f78=\", newText: \"423\u00f75=\" },\n    { oldText: \"644\u00f78=\", newText: \"983\u00f72=\" },\n    { oldText: \"380\u00f73=\", newText: \"687\u00f78=\" },\n    { oldText: \"934\u00f77=\", newText: \"962\u00f74=\" },\n    { oldText: \"476\u00f76=\", newText: \"999\u00f74=\" },\n    { oldText: \"542\u00f78=\", newText: \"205\u00f72=\" },\n    { oldText: \"312\u00f79=\", newText: \"611\u00f72=\" },\n    { oldText: \"323\u00f79=\", newText: \"254\u00f73=\" },\n    { oldText: \"494\u00f74=\", newText: \"548\u00f73=\" },\n    { oldText: \"272\u00f76=\", newText: \"758\u00f79=\" },\n    { oldText: \"838\u00f75=\", newText: \"510\u00f76=\" },\n    { oldText: \"208\u00f73=\", newText: \"772\u00f75=\" },\n    { oldText: \"662\u00f77=\", newText: \"730\u00f72=\" },\n    { oldText: \"981\u00f72=\", newText: \"724\u00f79=\" },\n    { oldText: \"826\u00f76=\", newText: \"324\u00f72=\" },\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each division-problem placeholder in the table cells with its\n# new value, per the commit diff (one Find/Replace per unique cell text).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"963\u00f77=\"; New = \"708\u00f76=\" },\n    @{ Old = \"496\u00f74=\"; New = \"555\u00f76=\" },\n    @{ Old = \"877\u00f79=\"; New = \"732\u00f75=\" },\n    @{ Old = \"708\u00f73=\"; New = \"257\u00f73=\" },\n    @{ Old = \"316\u00f74=\"; New = \"463\u00f72=\" },\n    @{ Old = \"314\u00f78=\"; New = \"558\u00f77=\" },\n    @{ Old = \"855\u00f74=\"; New = \"233\u00f73=\" },\n    @{ Old = \"703\u00f73=\"; New = \"705\u00f77=\" },\n    @{ Old = \"128\u00f75=\"; New = \"748\u00f72=\" },\n    @{ Old = \"634\u00f74=\"; New = \"694\u00f75=\" },\n    @{ Old = \"971\u00f78=\"; New = \"423\u00f75=\" },\n    @{ Old = \"644\u00f78=\"; New = \"983\u00f72=\" },\n    @{ Old = \"380\u00f73=\"; New = \"687\u00f78=\" },\n    @{ Old = \"934\u00f77=\"; New = \"962\u00f74=\" },\n    @{ Old = \"476\u00f76=\"; New = \"999\u00f74=\" },\n    @{ Old = \"542\u00f78=\"; New = \"205\u00f72=\" },\n    @{ Old = \"312\u00f79=\"; New = \"611\u00f72=\" },\n    @{ Old = \"323\u00f79=\"; New = \"254\u00f73=\" },\n    @{ Old = \"494\u00f74=\"; New = \"548\u00f73=\" },\n    @{ Old = \"272\u00f76=\"; New = \"758\u00f79=\" },\n    @{ Old = \"838\u00f75=\"; New = \"510\u00f76=\" },\n    @{ Old = \"208\u00f73=\"; New = \"772\u00f75=\" },\n    @{ Old = \"662\u00f77=\"; New = \"730\u00f72=\" },\n    @{ Old = \"981\u00f72=\"; New = \"724\u00f79=\" },\n    @{ Old = \"826\u00f76=\"; New = \"324\u00f72=\" },\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop -- each source string is unique, no wraparound needed\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n    if (-not $found) {\n        throw \"Text not found: $($r.Old)\"\n    }\n}\n"}
